$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "62.406.69"
$ws.Range("E2").Value = "  +2.14%  "
$ws.Range("D3").Value = "3.432.53"
$ws.Range("E3").Value = "  +3.31%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'406.91"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +2.71%  "
$ws.Range("D6").Value = "'130.91"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +5.41%  "
$ws.Range("D7").Value = "'0.600"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +2.98%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +0.07%  "
$ws.Range("E9").Value = "  +8.01%  "
$ws.Range("D10").Value = "'0.141"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +20.85%  "
$ws.Range("D11").Value = "'42.42"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +5.67%  "
$ws.Range("E12").Value = "  +0.66%  "
$ws.Range("D13").Value = "'8.61"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +5.93%  "
$ws.Range("D14").Value = "'19.96"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.17%  "
$ws.Range("D15").Value = "3.451.60"
$ws.Range("E15").Value = "  +2.90%  "
$ws.Range("D16").Value = "62.477.34"
$ws.Range("E16").Value = "  +2.52%  "
$ws.Range("D17").Value = "'11.59"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.91%  "
$ws.Range("B18").Value = "Polygon"
$ws.Range("C18").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D18").Value = "'1.02"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.34%  "
$ws.Range("B19").Value = "ShibaInu"
$ws.Range("C19").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D19").Value = "'0.0000167"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +36.62%  "
$ws.Range("E20").Value = "  +1.54%  "
$ws.Range("D21").Value = "'84.78"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +7.71%  "
$ws.Range("D22").Value = "'315.34"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +7.01%  "
$ws.Range("E23").Value = "  +2.83%  "
$ws.Range("D24").Value = "'3.19"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +2.31%  "
$ws.Range("E25").Value = "  +1.77%  "
$ws.Range("D26").Value = "'30.07"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.42%  "
$ws.Range("D27").Value = "'8.18"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +1.98%  "
$ws.Range("E28").Value = "  +6.55%  "
$ws.Range("D29").Value = "'2.71"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +9.69%  "
$ws.Range("E30").Value = "  +2.49%  "
$ws.Range("D31").Value = "'44.30"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +9.09%  "
$ws.Range("E32").Value = "  +3.53%  "
$ws.Range("E33").Value = "  +3.63%  "
$ws.Range("E34").Value = "  -0.08%  "
$ws.Range("E35").Value = "  +4.17%  "
$ws.Range("D36").Value = "'51.34"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -0.79%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +0.45%  "
$ws.Range("E38").Value = "  +4.66%  "
$ws.Range("D39").Value = "'3.35"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +0.93%  "
$ws.Range("E40").Value = "  +17.04%  "
$ws.Range("D41").Value = "'143.81"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +6.26%  "
$ws.Range("E42").Value = "  +4.86%  "
$ws.Range("E43").Value = "  +3.50%  "
$ws.Range("D44").Value = "'17.01"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +4.41%  "
$ws.Range("E45").Value = "  +3.93%  "
$ws.Range("E46").Value = "  +0.72%  "
$ws.Range("D47").Value = "'21.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +2.16%  "
$ws.Range("D48").Value = "2.109.13"
$ws.Range("E48").Value = "  +1.38%  "
$ws.Range("D49").Value = "'2.03"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +12.61%  "
$ws.Range("E50").Value = "  +2.26%  "
$ws.Range("E51").Value = "  +32.98%  "
